$wb = $excel.ActiveWorkbook

# --- Insert "Sheet" right after "Only US-Result sheet " ---
$after1 = $wb.Worksheets.Item(1)
$sheetNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after1)
$sheetNew.Name = "Sheet"

$sheetNew.Range("A1").Value = "total_time_3_5"
$sheetNew.Range("B1").Value = "avg_time_3_5"
$sheetNew.Range("C1").Value = "total_time_4"
$sheetNew.Range("D1").Value = "avg_time_4"
$sheetNew.Range("A1:D1").Font.Bold = $true
$sheetNew.Range("A1:D1").HorizontalAlignment = -4108
$sheetNew.Range("A1:D1").VerticalAlignment = -4160
$sheetNew.Range("A1:D1").Borders.LineStyle = 1

$sheetNew.Range("A2").Value = 1000
$sheetNew.Range("B2").Value = 1000
$sheetNew.Range("C2").Value = 1000
$sheetNew.Range("D2").Value = 1000

# --- Append "TimeConsupNoAnn" as the last sheet ---
$lastIdx = $wb.Worksheets.Count
$afterLast = $wb.Worksheets.Item($lastIdx)
$timeSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterLast)
$timeSheet.Name = "TimeConsupNoAnn"

$timeSheet.Range("A1").Value = "Version"
$timeSheet.Range("B1").Value = "Groups"
$timeSheet.Range("C1").Value = "Total Time"
$timeSheet.Range("D1").Value = "Avg Time"

$timeSheet.Range("A2").Value = "v3.5"
$timeSheet.Range("B2").Value = "All sets"
$timeSheet.Range("C2").Value = 1005
$timeSheet.Range("D2").Value = 1005

$timeSheet.Range("A3").Value = "v4"
$timeSheet.Range("B3").Value = "All sets"
$timeSheet.Range("C3").Value = 1000
$timeSheet.Range("D3").Value = 1000

$timeSheet.Range("A4").Value = "v3.5"
$timeSheet.Range("B4").Value = "00_#G03#"
$timeSheet.Range("C4").Value = 1005
$timeSheet.Range("D4").Value = 1005

$timeSheet.Range("A5").Value = "v4"
$timeSheet.Range("B5").Value = "00_#G03#"
$timeSheet.Range("C5").Value = 1000
$timeSheet.Range("D5").Value = 1000

$timeSheet.Columns.Item(1).ColumnWidth = 13.5
$timeSheet.Columns.Item(2).ColumnWidth = 12
$timeSheet.Columns.Item(3).ColumnWidth = 18
$timeSheet.Columns.Item(4).ColumnWidth = 15

$timeSheet.Range("A1:D1").Font.Bold = $true
$timeSheet.Range("A1:D1").HorizontalAlignment = -4108
$timeSheet.Range("A1:D1").VerticalAlignment = -4160
$timeSheet.Range("A1:D1").Borders.LineStyle = 1

$timeSheet.Range("A2:D5").Font.Bold = $true
$timeSheet.Range("A2:D5").WrapText = $true
$timeSheet.Range("A2:D5").HorizontalAlignment = -4131
$timeSheet.Range("A2:D5").VerticalAlignment = -4160
$timeSheet.Range("A2:D5").Borders.LineStyle = 1

$timeSheet.Range("A1:D1").AutoFilter() | Out-Null

$timeSheet.Range("A1").Select()
$timeSheet.Application.ActiveWindow.FreezePanes = $false
$timeSheet.Range("A2").Select()
$timeSheet.Application.ActiveWindow.FreezePanes = $true
$timeSheet.Range("A1").Select()
